$d = $word.ActiveDocument

$d.Content.Find.Execute("811×5=", $true, $true, $false, $false, $false, $true, 1, $false, "347×2=", 2) | Out-Null
$d.Content.Find.Execute("275×9=", $true, $true, $false, $false, $false, $true, 1, $false, "295×6=", 2) | Out-Null
$d.Content.Find.Execute("748×7=", $true, $true, $false, $false, $false, $true, 1, $false, "397×3=", 2) | Out-Null
$d.Content.Find.Execute("187×7=", $true, $true, $false, $false, $false, $true, 1, $false, "446×3=", 2) | Out-Null
$d.Content.Find.Execute("458×5=", $true, $true, $false, $false, $false, $true, 1, $false, "693×4=", 2) | Out-Null
$d.Content.Find.Execute("326×9=", $true, $true, $false, $false, $false, $true, 1, $false, "280×4=", 2) | Out-Null
$d.Content.Find.Execute("193×5=", $true, $true, $false, $false, $false, $true, 1, $false, "702×3=", 2) | Out-Null
$d.Content.Find.Execute("463×5=", $true, $true, $false, $false, $false, $true, 1, $false, "547×4=", 2) | Out-Null
$d.Content.Find.Execute("355×7=", $true, $true, $false, $false, $false, $true, 1, $false, "363×3=", 2) | Out-Null
$d.Content.Find.Execute("922×7=", $true, $true, $false, $false, $false, $true, 1, $false, "390×9=", 2) | Out-Null
$d.Content.Find.Execute("466×3=", $true, $true, $false, $false, $false, $true, 1, $false, "579×3=", 2) | Out-Null
$d.Content.Find.Execute("439×9=", $true, $true, $false, $false, $false, $true, 1, $false, "320×8=", 2) | Out-Null
$d.Content.Find.Execute("435×8=", $true, $true, $false, $false, $false, $true, 1, $false, "550×9=", 2) | Out-Null
$d.Content.Find.Execute("218×3=", $true, $true, $false, $false, $false, $true, 1, $false, "971×9=", 2) | Out-Null
$d.Content.Find.Execute("273×9=", $true, $true, $false, $false, $false, $true, 1, $false, "582×8=", 2) | Out-Null
$d.Content.Find.Execute("623×2=", $true, $true, $false, $false, $false, $true, 1, $false, "650×5=", 2) | Out-Null
$d.Content.Find.Execute("897×5=", $true, $true, $false, $false, $false, $true, 1, $false, "976×9=", 2) | Out-Null
$d.Content.Find.Execute("777×2=", $true, $true, $false, $false, $false, $true, 1, $false, "234×7=", 2) | Out-Null
$d.Content.Find.Execute("841×3=", $true, $true, $false, $false, $false, $true, 1, $false, "540×9=", 2) | Out-Null
$d.Content.Find.Execute("588×8=", $true, $true, $false, $false, $false, $true, 1, $false, "498×7=", 2) | Out-Null
$d.Content.Find.Execute("729×8=", $true, $true, $false, $false, $false, $true, 1, $false, "356×9=", 2) | Out-Null
$d.Content.Find.Execute("458×8=", $true, $true, $false, $false, $false, $true, 1, $false, "531×8=", 2) | Out-Null
$d.Content.Find.Execute("911×8=", $true, $true, $false, $false, $false, $true, 1, $false, "845×3=", 2) | Out-Null
$d.Content.Find.Execute("631×5=", $true, $true, $false, $false, $false, $true, 1, $false, "121×5=", 2) | Out-Null
$d.Content.Find.Execute("766×8=", $true, $true, $false, $false, $false, $true, 1, $false, "917×2=", 2) | Out-Null
